$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154; this shifts existing rows 154..238 down to 155..239
$ws.Rows(154).Insert()

# Populate the newly inserted row 154 with the new record.
# Columns A,B,C,E,F,G,H,I,J,R are identical to the surrounding rows (same market/product block).
$ws.Cells.Item(154, 1).Value = 10
$ws.Cells.Item(154, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(154, 3).Value = "La Araucanía"
$ws.Cells.Item(154, 4).Value = 44719
$ws.Cells.Item(154, 5).Value = 9
$ws.Cells.Item(154, 6).Value = "Fruta"
$ws.Cells.Item(154, 7).Value = 100103
$ws.Cells.Item(154, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(154, 9).Value = 100103002
$ws.Cells.Item(154, 10).Value = "Ciruela"
$ws.Cells.Item(154, 11).Value = "Pink Delight"
$ws.Cells.Item(154, 12).Value = "Primera"
$ws.Cells.Item(154, 13).Value = 65
$ws.Cells.Item(154, 14).Value = 12000
$ws.Cells.Item(154, 15).Value = 12000
$ws.Cells.Item(154, 16).Value = 12000
$ws.Cells.Item(154, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(154, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(154, 19).Value = 667
$ws.Cells.Item(154, 20).Value = 18
